# "Remove code as classes"
# The "Classes" sheet listed a number of TED/ePO code-list types (AccountFormat,
# AccountType, BuyerCategoryType, BuyerRoleType, Country, Currency, EOGroupCodeType,
# IdentifierProvider, LegalForm (twice), ProcedureChoiceJustificationCode,
# ProcedureType, ReservedContract, ReservedContractType, SubmissionLanguage and the
# generic "Code" entry) as if they were classes. They are in fact codes, not
# classes, so their rows are deleted from the sheet entirely. Deleting the rows
# also drops the shared-string entries that become unused, and every remaining
# shared-string based cell (row labels on "Classes" plus the TED-XSD comment
# column on "Triples") is renumbered accordingly once the workbook is saved.

$wb = $excel.ActiveWorkbook
$wsClasses = $wb.Worksheets.Item("Classes")
$wsTriples = $wb.Worksheets.Item("Triples")

# Rows on "Classes" whose B-column text is one of the removed "code" entries
# (includes both duplicate "LegalForm" rows, 31 and 32). Deleted from the
# bottom up so earlier row numbers in the list stay valid as we go.
$rowsToDelete = @(4,5,7,9,13,17,19,23,29,31,32,37,38,46,47,49) | Sort-Object -Descending

foreach ($r in $rowsToDelete) {
    $wsClasses.Rows.Item($r).Delete()
}

# Restore the view state: "Classes" becomes the active/selected sheet with B2
# selected, and the "Triples" sheet keeps a selection over the block that was
# being reviewed (B8:E15) without remaining the tab-selected sheet.
$wsTriples.Range("B8:E15").Select()
$wsClasses.Activate()
$wsClasses.Range("B2").Select()
